# Agil - Copia de Preguntas_Examen.xlsx  -- "Fichero de Preguntas - Correccion"
#
# This script:
#  1. Inserts a new question row right after the existing
#     "Which basic Agile quality practice reduces bottlenecks and ensures
#     consistency?" row (original row 337), with a different set of options
#     that tests "Peer-review and pairing".
#  2. Renumbers the "Nº" column (A) sequentially for every data row, since
#     the original sheet had a couple of pre-existing gaps.
#  3. Corrects the "Respuesta Correcta" (column D) for four questions whose
#     marked-correct answer was wrong.
#  4. Restores the view state (top-left cell / selection) to match the
#     saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert the new row.
#    Before the insert, "Which basic Agile quality practice reduces
#    bottlenecks and ensures consistency?" lives on row 337, immediately
#    followed by "Which statement is a value from the Agile Manifesto?" on
#    row 338. We insert a new blank row at 338, pushing everything down by
#    one, and fill it in.
# ---------------------------------------------------------------------
$anchor = $ws.Cells.Find("Which basic Agile quality practice reduces bottlenecks and ensures consistency?")
$newRow = $anchor.Row + 1

$ws.Rows.Item($newRow).Insert()

$ws.Cells.Item($newRow, 2).Value = "Which basic Agile quality practice reduces bottlenecks and ensures consistency?"
$ws.Cells.Item($newRow, 3).Value = "Definition of Done`nEstablish flow`nCollective owership and standards`nPeer-review and pairing"
$ws.Cells.Item($newRow, 4).Value = "Peer-review and pairing"
$ws.Cells.Item($newRow, 5).Value = 0
$ws.Cells.Item($newRow, 6).Value = 0

# ---------------------------------------------------------------------
# 2. Renumber column A (the "Nº" column) sequentially: 1, 2, 3, ...
#    for every data row (row 2 is the first data row, right after the
#    header row 1).
# ---------------------------------------------------------------------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = ($r - 1)
}

# ---------------------------------------------------------------------
# 3. Fix the four incorrect "Respuesta Correcta" entries.
# ---------------------------------------------------------------------
$fix1 = $ws.Cells.Find("A confidence vote is taken at the end of PI Planning after dependencies are resolved and risks are addressed. What best describes the process of the confidence vote?")
$ws.Cells.Item($fix1.Row, 4).Value = "The teams and the ARTs vote"

$fix2 = $ws.Cells.Find("Which statement applies to uncommitted objectives?")
$ws.Cells.Item($fix2.Row, 4).Value = "They are counted when calculating load"

$fix3 = $ws.Cells.Find("Why do Business Owners assign business value to team PI Objectives?")
$ws.Cells.Item($fix3.Row, 4).Value = "To empower teams to make decisions around work"

$fix4 = $ws.Cells.Find("What is one way to describe a cross-functional Agile Team?")
$ws.Cells.Item($fix4.Row, 4).Value = "They are optimized for communication and delivery of value"

# ---------------------------------------------------------------------
# 4. Restore the view state saved with the workbook.
# ---------------------------------------------------------------------
$ws.Range("A2:XFD381").Select()
$excel.ActiveWindow.ScrollRow = 349
